$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.549.48"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "3.082.68"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.49"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.14"
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.077.27"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.41"
$ws.Range("E11").Value = "  +2.69%  "
$ws.Range("E12").Value = "  -2.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.06"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000225"
$ws.Range("E14").Value = "  +3.61%  "
$ws.Range("D15").Value = "3.581.53"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "63.503.05"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "3.079.68"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.61"
$ws.Range("E20").Value = "  -2.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.50"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.702"
$ws.Range("E22").Value = "  -2.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.73"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.26"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.97"
$ws.Range("E28").Value = "  -5.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.26"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("E31").Value = "  -3.39%  "
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "57.96"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  -7.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.48"
$ws.Range("E35").Value = "  +6.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "494.57"
$ws.Range("E36").Value = "  -3.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.03"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "3.263.33"
$ws.Range("E38").Value = "  +3.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0405"
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0801"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.16"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.46"
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.63"
$ws.Range("E47").Value = "  +2.24%  "
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("D49").Value = "0.0₃0532"
$ws.Range("E49").Value = "  +5.89%  "
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("E51").Value = "  -0.24%  "
